# Atualizacao de bases das ligas - apply updated odds/results data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 86-87: swap match records (Podbeskidzie vs Gornik Leczna / Miedz Legnica vs Odra Opole) and refresh odds
$ws.Range("B86").Value = 6803738
$ws.Range("E86").Value = "Podbeskidzie Bielsko Biala"
$ws.Range("F86").Value = "Gornik Leczna"
$ws.Range("H86").Value = 1
$ws.Range("I86").Value = "D"
$ws.Range("L86").Value = 3.8
$ws.Range("M86").Value = 1.666
$ws.Range("N86").Value = 3.8
$ws.Range("O86").Value = 4.75
$ws.Range("P86").Value = -0.75
$ws.Range("Q86").Value = 1.825
$ws.Range("R86").Value = 1.975
$ws.Range("S86").Value = 2.5
$ws.Range("T86").Value = 1.825
$ws.Range("U86").Value = 1.975
$ws.Range("W86").Value = 2.8
$ws.Range("X86").Value = -1
$ws.Range("Z86").Value = 0.9750000000000001
$ws.Range("AA86").Value = -1
$ws.Range("AB86").Value = 0.9750000000000001
$ws.Range("B87").Value = 6803740
$ws.Range("E87").Value = "Miedz Legnica"
$ws.Range("F87").Value = "Odra Opole"
$ws.Range("H87").Value = 2
$ws.Range("I87").Value = "A"
$ws.Range("L87").Value = 3.75
$ws.Range("M87").Value = 1.909
$ws.Range("N87").Value = 3.5
$ws.Range("O87").Value = 3.5
$ws.Range("P87").Value = -0.5
$ws.Range("Q87").Value = 1.975
$ws.Range("R87").Value = 1.825
$ws.Range("W87").Value = -1
$ws.Range("X87").Value = 2.5
$ws.Range("Z87").Value = 0.825
$ws.Range("AA87").Value = 0.8999999999999999
$ws.Range("AB87").Value = -1

# Rows 168-169: swap match records (Wisla Krakow vs Gornik Leczna / Odra Opole vs Stal Rzeszow) and refresh odds
$ws.Range("B168").Value = 6803794
$ws.Range("E168").Value = "Wisla Krakow"
$ws.Range("F168").Value = "Gornik Leczna"
$ws.Range("G168").Value = 4
$ws.Range("H168").Value = 0
$ws.Range("I168").Value = "H"
$ws.Range("J168").Value = 1.4
$ws.Range("K168").Value = 4.75
$ws.Range("L168").Value = 7
$ws.Range("M168").Value = 1.363
$ws.Range("N168").Value = 4.75
$ws.Range("O168").Value = 7.5
$ws.Range("P168").Value = -1.25
$ws.Range("Q168").Value = 1.8
$ws.Range("R168").Value = 2
$ws.Range("S168").Value = 2.75
$ws.Range("T168").Value = 1.775
$ws.Range("U168").Value = 2.025
$ws.Range("V168").Value = 0.363
$ws.Range("W168").Value = -1
$ws.Range("Y168").Value = 0.8
$ws.Range("Z168").Value = -1
$ws.Range("AA168").Value = 0.7749999999999999
$ws.Range("AB168").Value = -1
$ws.Range("B169").Value = 6803793
$ws.Range("E169").Value = "Odra Opole"
$ws.Range("F169").Value = "Stal Rzeszow"
$ws.Range("G169").Value = 1
$ws.Range("H169").Value = 1
$ws.Range("I169").Value = "D"
$ws.Range("J169").Value = 2.05
$ws.Range("K169").Value = 3.4
$ws.Range("L169").Value = 3.5
$ws.Range("M169").Value = 2.1
$ws.Range("N169").Value = 3.4
$ws.Range("O169").Value = 3.4
$ws.Range("P169").Value = -0.25
$ws.Range("Q169").Value = 1.825
$ws.Range("R169").Value = 1.975
$ws.Range("S169").Value = 2.5
$ws.Range("T169").Value = 1.95
$ws.Range("U169").Value = 1.85
$ws.Range("V169").Value = -1
$ws.Range("W169").Value = 2.4
$ws.Range("Y169").Value = -0.5
$ws.Range("Z169").Value = 0.4875
$ws.Range("AA169").Value = -1
$ws.Range("AB169").Value = 0.8500000000000001

# Rows 257-258: swap match records (Stal Rzeszow vs Podbeskidzie / Polonia Warsaw vs Zaglebie Sosnowiec) and refresh odds
$ws.Range("B257").Value = 6803855
$ws.Range("E257").Value = "Stal Rzeszow"
$ws.Range("F257").Value = "Podbeskidzie Bielsko Biala"
$ws.Range("G257").Value = 2
$ws.Range("H257").Value = 2
$ws.Range("J257").Value = 2.25
$ws.Range("K257").Value = 3.4
$ws.Range("L257").Value = 2.9
$ws.Range("M257").Value = 2.4
$ws.Range("N257").Value = 3.25
$ws.Range("O257").Value = 2.7
$ws.Range("P257").Value = 0
$ws.Range("Q257").Value = 1.775
$ws.Range("R257").Value = 2.025
$ws.Range("T257").Value = 1.8
$ws.Range("U257").Value = 2
$ws.Range("W257").Value = 2.25
$ws.Range("Y257").Value = 0
$ws.Range("Z257").Value = 0
$ws.Range("AA257").Value = 0.8
$ws.Range("B258").Value = 6805658
$ws.Range("E258").Value = "Polonia Warsaw"
$ws.Range("F258").Value = "Zaglebie Sosnowiec"
$ws.Range("G258").Value = 3
$ws.Range("H258").Value = 3
$ws.Range("J258").Value = 1.95
$ws.Range("K258").Value = 3.5
$ws.Range("L258").Value = 3.5
$ws.Range("M258").Value = 1.909
$ws.Range("N258").Value = 3.5
$ws.Range("O258").Value = 3.6
$ws.Range("P258").Value = -0.5
$ws.Range("Q258").Value = 1.975
$ws.Range("R258").Value = 1.825
$ws.Range("T258").Value = 1.925
$ws.Range("U258").Value = 1.875
$ws.Range("W258").Value = 2.5
$ws.Range("Y258").Value = -1
$ws.Range("Z258").Value = 0.825
$ws.Range("AA258").Value = 0.925

# Rows 279-280: swap match records (Polonia Warsaw vs GKS Katowice / Miedz Legnica vs Znicz Pruszkw) and refresh odds
$ws.Range("B279").Value = 6884022
$ws.Range("E279").Value = "Polonia Warsaw"
$ws.Range("F279").Value = "GKS Katowice"
$ws.Range("J279").Value = 2.5
$ws.Range("K279").Value = 3.4
$ws.Range("L279").Value = 2.5
$ws.Range("M279").Value = 2.8
$ws.Range("O279").Value = 2.25
$ws.Range("P279").Value = 0.25
$ws.Range("Q279").Value = 1.825
$ws.Range("R279").Value = 2.025
$ws.Range("S279").Value = 2.5
$ws.Range("T279").Value = 2
$ws.Range("U279").Value = 1.85
$ws.Range("X279").Value = 1.25
$ws.Range("Z279").Value = 1.025
$ws.Range("AA279").Value = 1
$ws.Range("B280").Value = 6803960
$ws.Range("E280").Value = "Miedz Legnica"
$ws.Range("F280").Value = "Znicz Pruszkw"
$ws.Range("J280").Value = 1.85
$ws.Range("K280").Value = 3.5
$ws.Range("L280").Value = 4
$ws.Range("M280").Value = 1.909
$ws.Range("O280").Value = 3.8
$ws.Range("P280").Value = -0.5
$ws.Range("Q280").Value = 1.925
$ws.Range("R280").Value = 1.875
$ws.Range("S280").Value = 2.25
$ws.Range("T280").Value = 1.8
$ws.Range("U280").Value = 2
$ws.Range("X280").Value = 2.8
$ws.Range("Z280").Value = 0.875
$ws.Range("AA280").Value = 0.8

# Odds corrections for individual fixtures
$ws.Range("T284").Value = 2
$ws.Range("U284").Value = 1.85
$ws.Range("M288").Value = 1.571
$ws.Range("N288").Value = 3.8
$ws.Range("O288").Value = 4.75
$ws.Range("P288").Value = -1
$ws.Range("Q288").Value = 2.05
$ws.Range("R288").Value = 1.8
$ws.Range("M292").Value = 4.333
$ws.Range("N292").Value = 4.2
$ws.Range("O292").Value = 1.533
$ws.Range("Q292").Value = 1.875
$ws.Range("R292").Value = 1.975
